{"js": "// Remove the placeholder run \"vnpt.SiteAddress\" that follows the\n// \"\u0110\u1ecba ch\u1ec9: \" label in the \"B\u00ean A\" address line, leaving the label\n// paragraph intact but without the merge-field placeholder text.\nconst results = context.document.body.search(\"vnpt.SiteAddress\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the placeholder run \"vnpt.SiteAddress\" that follows the\n# \"\u0110\u1ecba ch\u1ec9: \" label in the \"B\u00ean A\" address line, leaving the label\n# paragraph intact but without the merge-field placeholder text.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"vnpt.SiteAddress\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1\n\nwhile ($find.Execute()) {\n    if (-not $find.Found) { break }\n    $find.Parent.Delete()\n}\n"}
